$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Password test-data cell (B2): was a bare numeric placeholder (12345),
# now the actual forgot-password PIN value used by the new automation flow.
$ws.Range("B2").Value = "Bhagwan@123"

# Excel auto-links "word@word"-shaped text typed into a cell; recreate that
# mailto hyperlink (mirrors the existing A2 email hyperlink) and restore the
# shared Hyperlink cell style that Add() nudges off A2's style.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Bhagwan@123", "", "", "Bhagwan@123")
$ws.Range("B2").Style = $ws.Range("A2").Style

# The active selection moved from A2 to D2 (the time_from input) while fixing
# the input-time error mentioned in the commit message.
$ws.Range("D2").Select() | Out-Null
